$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 10
$ws.Range("H2").Value = "Discrete"
